# Updated C3DC Regression and Smoke suites
# The "TreatmentTab" query text (in cell B5) had a redundant CONCAT() wrapper
# around the REPLACE(...) call for the "Treatment Agent" column. Remove the
# redundant CONCAT() wrapper, leaving a straightforward REPLACE(...) call.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$treatmentCell = $ws.Range("B5")
$originalQuery = $treatmentCell.Value2

$oldFragment = "CONCAT(REPLACE(trt.treatment_agent, ';', ', '))"
$newFragment = "REPLACE(trt.treatment_agent, ';', ', ')"

$updatedQuery = $originalQuery.Replace($oldFragment, $newFragment)
$treatmentCell.Value2 = $updatedQuery

# Re-apply (re-assert) the font on the corrected cell; this mirrors the
# author re-touching the formatting on B5 after editing its text, which is
# why B5 ends up with its own distinct style entry separate from the other
# query cells (B6, B7) that kept their original, shared style.
$treatmentCell.Font.Name = "Calibri"
$treatmentCell.Font.Size = 12
$treatmentCell.WrapText = $true

# Update the selection / active cell to B5, matching the saved view state.
$treatmentCell.Select()
